$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.074.77'
$ws.Range("E2").Value = '  -1.29%  '

$ws.Range("D3").Value = '2.465.96'
$ws.Range("E3").Value = '  -2.93%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.47'
$ws.Range("E5").Value = '  -1.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.75'
$ws.Range("E6").Value = '  -2.90%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D9").Value = '2.465.96'
$ws.Range("E9").Value = '  -2.89%  '

$ws.Range("E10").Value = '  -3.00%  '

$ws.Range("E11").Value = '  -0.34%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.93'
$ws.Range("E12").Value = '  -3.38%  '

$ws.Range("E13").Value = '  -4.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.57'
$ws.Range("E14").Value = '  -3.75%  '

$ws.Range("D15").Value = '2.917.07'
$ws.Range("E15").Value = '  -1.91%  '

$ws.Range("D16").Value = '66.910.29'
$ws.Range("E16").Value = '  -1.16%  '

$ws.Range("E17").Value = '  -4.95%  '

$ws.Range("D18").Value = '2.427.03'
$ws.Range("E18").Value = '  -3.91%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.10'
$ws.Range("E19").Value = '  -5.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.56'
$ws.Range("E20").Value = '  -5.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '353.57'
$ws.Range("E21").Value = '  -4.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.03'
$ws.Range("E22").Value = '  -3.14%  '

$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.99'
$ws.Range("E24").Value = '  -3.66%  '

$ws.Range("E25").Value = '  -7.63%  '

$ws.Range("E26").Value = '  -7.24%  '

$ws.Range("E27").Value = '  -7.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.996'
$ws.Range("E28").Value = '  -60.12%  '

$ws.Range("D29").Value = '2.580.25'
$ws.Range("E29").Value = '  -3.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '515.74'
$ws.Range("E30").Value = '  -5.00%  '

$ws.Range("D31").Value = '0.0₃0902'
$ws.Range("E31").Value = '  -7.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.75'
$ws.Range("E32").Value = '  -8.64%  '

$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.24'
$ws.Range("E33").Value = '  -6.52%  '

$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.77'
$ws.Range("E34").Value = '  -5.33%  '

$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("E36").Value = '  -9.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '158.61'
$ws.Range("E37").Value = '  -0.56%  '

$ws.Range("E38").Value = '  +0.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.34'
$ws.Range("E39").Value = '  -4.43%  '

$ws.Range("E40").Value = '  -6.47%  '

$ws.Range("E41").Value = '  -0.32%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.81'
$ws.Range("E42").Value = '  -7.02%  '

$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.67'
$ws.Range("E43").Value = '  -6.50%  '

$ws.Range("E44").Value = '  -7.30%  '

$ws.Range("E45").Value = '  -7.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.78'
$ws.Range("E46").Value = '  -1.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.43'
$ws.Range("E47").Value = '  -5.16%  '

$ws.Range("E48").Value = '  -7.21%  '

$ws.Range("E49").Value = '  -7.18%  '

$ws.Range("E50").Value = '  -12.68%  '

$ws.Range("E51").Value = '  -7.51%  '
